$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to be inserted at the top (before existing row 2, which holds the
# first data row). This shifts all existing data rows down by 5.
$topRows = @(
    @(-2.194540023803711, 3.003044128417969, 1.832991361618042),
    @(-2.376946449279785, 2.956967353820801, 2.026212930679321),
    @(-2.072798252105713, 2.900663375854492, 1.98937726020813),
    @(-2.095717430114746, 2.981966018676758, 1.861132502555847),
    @(-1.865831851959228, 2.951019287109375, 1.621297121047974)
)

# Insert 5 new blank rows above row 2, then fill them with the new data.
$insertRange = $ws.Range("A2:C6")
$insertRange.EntireRow.Insert()
# Row-insert copies formatting from the row above (the bold/bordered header),
# which the source data rows never had - strip it back to unstyled.
$insertRange.ClearFormats()

for ($i = 0; $i -lt $topRows.Length; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $topRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $topRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $topRows[$i][2]
}

# New rows appended at the bottom (after the last existing data row).
$bottomRows = @(
    @(-0.0590333938598632, 3.949368000030518, 1.155394554138184),
    @(-0.8675603866577148, 3.302557468414306, 1.166132688522339),
    @(-0.5926990509033203, 3.298477172851562, 1.354385137557983),
    @(-0.5156621932983398, 3.70509934425354, 1.191632270812989),
    @(-0.767481803894043, 3.210179328918457, 1.142881870269775)
)

$startRow = 27
for ($i = 0; $i -lt $bottomRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $bottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $bottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $bottomRows[$i][2]
}
